# BGDP BAU GDP workbook update: OECD "Data" sheet now reports Canada (CAN)
# GDP long-term forecast values instead of USA. Sharedstrings "USA"/"Value"
# are replaced by "CAN"/"VALUE", and the forecast numbers in column G are
# replaced with the new Canadian series (also formatted as whole numbers).
# The "BGDP" sheet's formulas reference 'OECD Data'!G2:G48 and recompute
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OECD Data")

# Header text for column G changes from "Value" to "VALUE".
$ws.Range("G1").Value = "VALUE"

# Location column changes from "USA" to "CAN" for every data row.
$ws.Range("A2:A48").Value = "CAN"

# New GDP long-term forecast values (Canada), one per year 2014-2060.
$newValues = @(
  1505495.47140106,
  1520564.3540934401,
  1542066.2592517899,
  1588338.69301313,
  1621635.5389815201,
  1657928.8942838099,
  1684775,
  1710503,
  1736864,
  1764255,
  1792774,
  1822387,
  1852965,
  1884677,
  1917532,
  1951559,
  1986793,
  2023273,
  2060808,
  2099401,
  2139034,
  2179688,
  2221353,
  2264030,
  2307720,
  2352425,
  2398151,
  2444905,
  2492690,
  2541479,
  2591227,
  2641896,
  2693479,
  2745988,
  2799434,
  2853831,
  2909216,
  2965657,
  3023227,
  3081985,
  3141981,
  3203277,
  3265966,
  3330145,
  3395881,
  3463205,
  3532133
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
  $row = 2 + $i
  $cell = $ws.Cells.Item($row, 7)
  $cell.Value = $newValues[$i]
  $cell.NumberFormat = "0"
}

# Leave the selection/active-sheet state the way the consultant left it:
# "About" selected C12, "OECD Data" had the refreshed column selected, and
# the final view rests on "BGDP" (now the active tab) with O23 selected.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("C12").Select()

$ws.Activate()
$ws.Range("G2:G48").Select()

$wsBgdp = $wb.Worksheets.Item("BGDP")
$wsBgdp.Activate()
$wsBgdp.Range("O23").Select()

$wb.Save()
